# edit.ps1 — apply the "s_plan" commit:
#   * update review-schedule values for rows 4-7 on sheet "杨瀚森"
#   * append a new row (row 8) for "高中单词8"
#   * add a new sheet "尹嘉禾" (after "杨瀚森") with a single data row for "21天list1"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update existing rows whose review dates / counts shifted.
#    Columns: A=index  B=start date  C=content  D=review count
#             E=last review  F=next review  G=deadline
# ---------------------------------------------------------------------
$updates = @(
    @{ Row = 4; D = 3; E = 43079; F = 43086; G = 43101 },
    @{ Row = 5; D = 3; E = 43081; F = 43088; G = 43103 },
    @{ Row = 6; D = 2; E = 43079; F = 43083; G = 43090 },
    @{ Row = 7; D = 2; E = 43081; F = 43085; G = 43092 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws1.Cells.Item($r, 4).Value = $u.D
    $ws1.Cells.Item($r, 5).Value = $u.E
    $ws1.Cells.Item($r, 6).Value = $u.F
    $ws1.Cells.Item($r, 7).Value = $u.G
}

# ---------------------------------------------------------------------
# 2. Append the new "高中单词8" row (row 8), copying the formatting of
#    row 7 first so the new cells pick up the same styles (index/date
#    number formats, borders, etc.) before the values are overwritten.
# ---------------------------------------------------------------------
$ws1.Range("A7:G7").Copy($ws1.Range("A8:G8"))

$ws1.Cells.Item(8, 1).Value = ""
$ws1.Cells.Item(8, 2).Value = 43081
$ws1.Cells.Item(8, 3).Value = "高中单词8"
$ws1.Cells.Item(8, 4).Value = 0
$ws1.Cells.Item(8, 5).Value = 43081
$ws1.Cells.Item(8, 6).Value = 43082
$ws1.Cells.Item(8, 7).Value = 43084

# ---------------------------------------------------------------------
# 3. Add the new sheet "尹嘉禾" right after "杨瀚森" and populate its
#    header row + the single "21天list1" data row (same shape/styles as
#    sheet 1's header/data rows).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "尹嘉禾"

$ws1.Range("B1:G1").Copy($ws2.Range("B1:G1"))
$ws1.Range("A8:G8").Copy($ws2.Range("A2:G2"))

$ws2.Cells.Item(2, 3).Value = "21天list1"

$wb.Worksheets.Item(1).Activate()
